$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "Testcases"

# Insert two new columns at the front (A and B) - existing columns A-D shift to C-F
$ws.Range("A1:B1").EntireColumn.Insert()

# New header row (row 1)
$ws.Range("A1").Value = "TFName"
$ws.Range("B1").Value = "TFDescription"

# New data rows
$ws.Range("A2").Value = "Login 1"
$ws.Range("A3").Value = "Login 2"
$ws.Range("B2").Value = "Login with Username and Password"
$ws.Range("B3").Value = "Login with Username and Password"

# Bold / colored font for the new TFDescription header cell
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").Font.Size = 10
$ws.Range("B1").Font.Name = "Menlo"
$ws.Range("B1").Font.Color = 8421376

# Re-point the hyperlinks (previously on A2/A3) to their new location C2/C3
$ws.Range("A2").Hyperlinks.Delete() | Out-Null
$ws.Range("A3").Hyperlinks.Delete() | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://app-eu.earthsquad.global/api/rest-auth/login/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://app-eu.earthsquad.global/api/rest-auth/login/") | Out-Null

# Restore the original "Link" cell style on the hyperlink cells (Hyperlinks.Add applies its own style)
$ws.Range("C2:C3").Style = "Link"

# Update the active selection to match the new layout
$ws.Range("E3").Select() | Out-Null
